$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the delimiter used in generated header names from '~' to '.'
# and collapse per-index vibrational wavenumber headers down to a single
# shared "vib_wavenumber" label (matching upstream commit "Changed default
# delimiter character from '~' to '.'").
$ws.Range("D1").Value = "elements.C"
$ws.Range("E1").Value = "elements.H"
$ws.Range("I1:BJ1").Value = "vib_wavenumber"

# Update the view/selection state to match the refreshed workbook.
$ws.Range("T1").Select()
